$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95..208 down to 96..209.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new record.
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = 'Vega Modelo de Temuco'
$ws.Range("C95").Value = 'La Araucanía'
$ws.Range("D95").Value2 = 44483
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = 100114013
$ws.Range("G95").Value = 'Zanahoria'
$ws.Range("H95").Value = 'Sin especificar'
$ws.Range("I95").Value = 'Primera'
$ws.Range("J95").Value = 50
$ws.Range("K95").Value = 9000
$ws.Range("L95").Value = 9000
$ws.Range("M95").Value = 9000
$ws.Range("N95").Value = '$/saco 20 kilos'
$ws.Range("O95").Value = 'Provincia del Elquí'
$ws.Range("P95").Value = 450
$ws.Range("Q95").Value = 20
$ws.Range("R95").Value = 'Hortaliza'
